# Update the "F" column (view/read counters) on the sheets that contain
# them, matching the regenerated data snapshot committed at 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 436
$ws.Range("F5").Value  = 8822
$ws.Range("F8").Value  = 657
$ws.Range("F9").Value  = 309
$ws.Range("F10").Value = 164
$ws.Range("F12").Value = 27
$ws.Range("F13").Value = 3689
$ws.Range("F14").Value = 51
$ws.Range("F17").Value = 3214
$ws.Range("F19").Value = 1127
$ws.Range("F21").Value = 219
$ws.Range("F22").Value = 2478
$ws.Range("F23").Value = 85

# Sheet "演出" (performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 36

# Sheet "全部类型" (all types combined)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 436
$ws.Range("F5").Value  = 8822
$ws.Range("F8").Value  = 657
$ws.Range("F9").Value  = 309
$ws.Range("F10").Value = 164
$ws.Range("F12").Value = 27
$ws.Range("F13").Value = 3689
$ws.Range("F14").Value = 51
$ws.Range("F17").Value = 3214
$ws.Range("F19").Value = 1127
$ws.Range("F21").Value = 219
$ws.Range("F22").Value = 2478
$ws.Range("F23").Value = 36
$ws.Range("F24").Value = 85
